# Rename the "Documentation" tab to "Help" and the
# "DesireEditedSequences" tab to "TargetedSearch".
$wb = $excel.ActiveWorkbook

$helpSheet = $wb.Worksheets.Item("Documentation")
$helpSheet.Name = "Help"

$tsearchSheet = $wb.Worksheets.Item("DesireEditedSequences")
$tsearchSheet.Name = "TargetedSearch"

# Update the selection on the Layout sheet from the whole column C
# to just cell C2 (this also happens to activate Layout briefly,
# which is corrected below by re-activating Help as the last step).
$layoutSheet = $wb.Worksheets.Item("Layout")
[void]$layoutSheet.Range("C2").Select()

# Make "Help" (previously "Documentation") the selected/active tab
# instead of "Amplicon", and drop its previous frozen/scrolled
# topLeftCell so the sheet opens back at A1 with A3:C3 selected.
$helpSheet.Activate()
[void]$helpSheet.Range("A3:C3").Select()
